$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 7 de Agosto de 2020 a las 16:33'

$ws.Cells.Item(4, 2).Value = 5036881
$ws.Cells.Item(4, 3).Value = 4702
$ws.Cells.Item(4, 4).Value = 2577938
$ws.Cells.Item(4, 5).Value = 2296070
$ws.Cells.Item(4, 7).Value = 69
$ws.Cells.Item(4, 8).Value = 162873

$ws.Cells.Item(6, 2).Value = 2046141
$ws.Cells.Item(6, 3).Value = 20732
$ws.Cells.Item(6, 4).Value = 1392940
$ws.Cells.Item(6, 5).Value = 611289
$ws.Cells.Item(6, 7).Value = 274
$ws.Cells.Item(6, 8).Value = 41912

$ws.Cells.Item(22, 2).Value = 215604
$ws.Cells.Item(22, 3).Value = 394
$ws.Cells.Item(22, 5).Value = 10152

$ws.Cells.Item(38, 2).Value = 77709
$ws.Cells.Item(38, 3).Value = 1173
$ws.Cells.Item(38, 4).Value = 41393
$ws.Cells.Item(38, 5).Value = 35057
$ws.Cells.Item(38, 7).Value = 13
$ws.Cells.Item(38, 8).Value = 1259

$ws.Cells.Item(42, 2).Value = 68614
$ws.Cells.Item(42, 3).Value = 111
$ws.Cells.Item(42, 4).Value = 64200
$ws.Cells.Item(42, 5).Value = 3831
$ws.Cells.Item(42, 7).Value = 3
$ws.Cells.Item(42, 8).Value = 583

$ws.Cells.Item(46, 1).Value = 'Guatemala'
$ws.Cells.Item(46, 2).Value = 55270
$ws.Cells.Item(46, 3).Value = 931
$ws.Cells.Item(46, 4).Value = 43135
$ws.Cells.Item(46, 5).Value = 9967
$ws.Cells.Item(46, 7).Value = 49
$ws.Cells.Item(46, 8).Value = 2168

$ws.Cells.Item(47, 1).Value = 'Singapur'
$ws.Cells.Item(47, 2).Value = 54797
$ws.Cells.Item(47, 3).Value = 242
$ws.Cells.Item(47, 4).Value = 48031
$ws.Cells.Item(47, 5).Value = 6739
$ws.Cells.Item(47, 8).Value = 27

$ws.Cells.Item(48, 2).Value = 52351
$ws.Cells.Item(48, 3).Value = 290
$ws.Cells.Item(48, 4).Value = 38087
$ws.Cells.Item(48, 5).Value = 12518
$ws.Cells.Item(48, 7).Value = 3
$ws.Cells.Item(48, 8).Value = 1746

$ws.Cells.Item(63, 2).Value = 27608
$ws.Cells.Item(63, 3).Value = 276
$ws.Cells.Item(63, 5).Value = 12935
$ws.Cells.Item(63, 7).Value = 5
$ws.Cells.Item(63, 8).Value = 626

$ws.Cells.Item(64, 2).Value = 26990
$ws.Cells.Item(64, 3).Value = 362
$ws.Cells.Item(64, 4).Value = 18918
$ws.Cells.Item(64, 5).Value = 7237
$ws.Cells.Item(64, 7).Value = 7
$ws.Cells.Item(64, 8).Value = 835

$ws.Cells.Item(86, 2).Value = 9503
$ws.Cells.Item(86, 3).Value = 35
$ws.Cells.Item(86, 5).Value = 390

$ws.Cells.Item(91, 2).Value = 7706
$ws.Cells.Item(91, 3).Value = 41
$ws.Cells.Item(91, 4).Value = 6484

$ws.Cells.Item(119, 4).Value = 2564
$ws.Cells.Item(119, 5).Value = 264

$ws.Cells.Item(120, 1).Value = 'Namibia'
$ws.Cells.Item(120, 2).Value = 2802
$ws.Cells.Item(120, 3).Value = 150
$ws.Cells.Item(120, 4).Value = 575
$ws.Cells.Item(120, 5).Value = 2211
$ws.Cells.Item(120, 7).Value = 1
$ws.Cells.Item(120, 8).Value = 16

$ws.Cells.Item(121, 1).Value = 'Cuba'
$ws.Cells.Item(121, 2).Value = 2775
$ws.Cells.Item(121, 4).Value = 2409
$ws.Cells.Item(121, 5).Value = 278
$ws.Cells.Item(121, 8).Value = 88

$ws.Cells.Item(122, 1).Value = 'Cabo Verde'
$ws.Cells.Item(122, 2).Value = 2734
$ws.Cells.Item(122, 4).Value = 2010
$ws.Cells.Item(122, 5).Value = 697
$ws.Cells.Item(122, 8).Value = 27

$ws.Cells.Item(142, 1).Value = 'Uganda'
$ws.Cells.Item(142, 2).Value = 1254
$ws.Cells.Item(142, 3).Value = 31
$ws.Cells.Item(142, 4).Value = 1113
$ws.Cells.Item(142, 5).Value = 135
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = 6

$ws.Cells.Item(143, 1).Value = 'Jordania'
$ws.Cells.Item(143, 2).Value = 1232
$ws.Cells.Item(143, 4).Value = 1171
$ws.Cells.Item(143, 5).Value = 50
$ws.Cells.Item(143, 8).Value = 11

$ws.Cells.Item(144, 1).Value = 'Liberia'
$ws.Cells.Item(144, 2).Value = 1224
$ws.Cells.Item(144, 4).Value = 705
$ws.Cells.Item(144, 5).Value = 441
$ws.Cells.Item(144, 8).Value = 78

$ws.Cells.Item(169, 4).Value = 309
$ws.Cells.Item(169, 5).Value = 44

$ws.Cells.Item(173, 1).Value = 'Islas Feroe'
$ws.Cells.Item(173, 2).Value = 291
$ws.Cells.Item(173, 3).Value = 25
$ws.Cells.Item(173, 4).Value = 192
$ws.Cells.Item(173, 5).Value = 99

$ws.Cells.Item(174, 1).Value = 'Eritrea'
$ws.Cells.Item(174, 2).Value = 282
$ws.Cells.Item(174, 4).Value = 225
$ws.Cells.Item(174, 5).Value = 57
$ws.Cells.Item(174, 8).Value = 0

$ws.Cells.Item(175, 1).Value = 'Guadalupe'
$ws.Cells.Item(175, 2).Value = 279
$ws.Cells.Item(175, 4).Value = 179
$ws.Cells.Item(175, 5).Value = 86
$ws.Cells.Item(175, 8).Value = 14

$ws.Cells.Item(176, 1).Value = 'Martinica'
$ws.Cells.Item(176, 2).Value = 276
$ws.Cells.Item(176, 4).Value = 98
$ws.Cells.Item(176, 5).Value = 163
$ws.Cells.Item(176, 8).Value = 15

$ws.Cells.Item(182, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(182, 2).Value = 188
$ws.Cells.Item(182, 3).Value = 25
$ws.Cells.Item(182, 4).Value = 53
$ws.Cells.Item(182, 5).Value = 132
$ws.Cells.Item(182, 8).Value = 3

$ws.Cells.Item(183, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(183, 2).Value = 176
$ws.Cells.Item(183, 3).Value = 16
$ws.Cells.Item(183, 4).Value = 86
$ws.Cells.Item(183, 5).Value = 74
$ws.Cells.Item(183, 8).Value = 16

$ws.Cells.Item(202, 1).Value = 'Santa Lucia'

$ws.Cells.Item(203, 1).Value = 'Timor Oriental'

Write-Host "Applied country/case updates"
